# fix bug upload file
# Adds a "Status" column (Aktif / Non-Aktif) to the upload-format template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data for column M ("Status")
$ws.Range("M1").Value = "Status"
$ws.Range("M2").Value = "Aktif"
$ws.Range("M3").Value = "Non-Aktif"

# Match the bordered header/body look already used by the other columns
# (K:L) by copying their formatting onto the new column, cell by cell so
# the existing table layout (K1:L3) is left untouched.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null

$ws.Range("L2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null

$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null

# Restore the values PasteSpecial(formats) leaves alone anyway, and make
# sure the text survived the paste (values already set above).
$ws.Range("M1").Value = "Status"
$ws.Range("M2").Value = "Aktif"
$ws.Range("M3").Value = "Non-Aktif"

$excel.CutCopyMode = $false

# Move the active selection the way the author's session ended up.
$ws.Range("I13").Select()
